$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update timestamps in F2:F27 of the "data" sheet ---
$newTimes = @(
  "2021-10-05 14:20:10.011698",
  "2021-10-05 14:20:10.011706",
  "2021-10-05 14:20:10.011709",
  "2021-10-05 14:20:10.011712",
  "2021-10-05 14:20:10.011715",
  "2021-10-05 14:20:10.011717",
  "2021-10-05 14:20:10.011720",
  "2021-10-05 14:20:10.011722",
  "2021-10-05 14:20:10.011725",
  "2021-10-05 14:20:10.011728",
  "2021-10-05 14:20:10.011730",
  "2021-10-05 14:20:10.011733",
  "2021-10-05 14:20:10.011735",
  "2021-10-05 14:20:10.011738",
  "2021-10-05 14:20:10.011740",
  "2021-10-05 14:20:10.011743",
  "2021-10-05 14:20:10.011745",
  "2021-10-05 14:20:10.011748",
  "2021-10-05 14:20:10.011751",
  "2021-10-05 14:20:10.011753",
  "2021-10-05 14:20:10.011756",
  "2021-10-05 14:20:10.011758",
  "2021-10-05 14:20:10.011761",
  "2021-10-05 14:20:10.011764",
  "2021-10-05 14:20:10.011766",
  "2021-10-05 14:20:10.011769"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add new "metadata" worksheet, placed after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Header row (bold + border style, matching "data" sheet header look)
$meta.Cells.Item(1,2).Value = "data_name"
$meta.Cells.Item(1,3).Value = "data_id"
$meta.Cells.Item(1,4).Value = "data_version"
$meta.Cells.Item(1,5).Value = "data_version_created"
$meta.Cells.Item(1,6).Value = "panel_query_time"
$meta.Cells.Item(1,7).Value = "panel_get_request"

$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row 2
$meta.Cells.Item(2,1).Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$meta.Cells.Item(2,2).Value = "Extreme early-onset hypertension"
$meta.Cells.Item(2,3).Value = 314
$meta.Cells.Item(2,4).Value = "'1.14"
$meta.Cells.Item(2,5).Value = "2020-11-13T13:35:47.245110Z"
$meta.Cells.Item(2,6).Value = "2021-10-05 14:20:10.007727"
$meta.Cells.Item(2,7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/314/?format=json"

# Restore "data" as the active sheet (matches unchanged bookViews/activeTab in the diff)
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null

Write-Host "Edit complete"
